$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 76; this shifts the existing
# rows 76-151 down to 78-153 (and keeps formatting of the row below).
$ws.Rows("76:77").Insert()

# New row 76: Femacal de La Calera / Coquimbo, Arandano (blue), Primera,
# Provincia de Quillota, $/bandeja 2 kilos
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44512
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100101
$ws.Range("H76").Value = "Berries"
$ws.Range("I76").Value = 100101001
$ws.Range("J76").Value = "Arándano (blue)"
$ws.Range("K76").Value = "Sin especificar"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 45
$ws.Range("N76").Value = 10000
$ws.Range("O76").Value = 10000
$ws.Range("P76").Value = 10000
$ws.Range("Q76").Value = "$/bandeja 2 kilos"
$ws.Range("R76").Value = "Provincia de Quillota"
$ws.Range("S76").Value = 5000
$ws.Range("T76").Value = 2

# New row 77: Femacal de La Calera / Coquimbo, Arandano (blue), Segunda,
# Provincia de Quillota, $/bandeja 2 kilos
$ws.Range("A77").Value = 3
$ws.Range("B77").Value = "Femacal de La Calera"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44512
$ws.Range("E77").Value = 5
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100101
$ws.Range("H77").Value = "Berries"
$ws.Range("I77").Value = 100101001
$ws.Range("J77").Value = "Arándano (blue)"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Segunda"
$ws.Range("M77").Value = 40
$ws.Range("N77").Value = 8000
$ws.Range("O77").Value = 8000
$ws.Range("P77").Value = 8000
$ws.Range("Q77").Value = "$/bandeja 2 kilos"
$ws.Range("R77").Value = "Provincia de Quillota"
$ws.Range("S77").Value = 4000
$ws.Range("T77").Value = 2
